# Update "想去人数" (number of people interested) counts that changed
# between scrapes, on both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> list of (cell, newValue) pairs to update.
$updates = @{
    "展览"   = @(
        @{ Cell = "F12"; Value = 1173 },
        @{ Cell = "F21"; Value = 7649 },
        @{ Cell = "F26"; Value = 2153 },
        @{ Cell = "F29"; Value = 185 },
        @{ Cell = "F35"; Value = 1780 },
        @{ Cell = "F42"; Value = 1882 }
    )
    "全部类型" = @(
        @{ Cell = "F14"; Value = 1173 },
        @{ Cell = "F24"; Value = 7649 },
        @{ Cell = "F29"; Value = 2153 },
        @{ Cell = "F32"; Value = 185 },
        @{ Cell = "F38"; Value = 1780 },
        @{ Cell = "F45"; Value = 1882 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Range($entry.Cell).Value = $entry.Value
    }
}
